# Apply the "Add data for 2022-06-20" update:
#  - Rename sheet / update title text from "...06-11" to "...06-12"
#  - Update the "June (through 06-11)" label to "June (through 06-12)"
#  - Bump the June row (row 7) and Total row (row 8) counts for columns B..I

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (this also updates the <sheet name="..."/> entry
# in workbook.xml).
$ws.Name = "Through 2022-06-12"

# Update the shared string used as the row label for June.
$ws.Range("A7").Value = "June (through 06-12)"

# Update June row (row 7) values for years 2015-2022 (columns B-I).
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 16
$ws.Range("D7").Value = 21
$ws.Range("E7").Value = 24
$ws.Range("F7").Value = 14
$ws.Range("G7").Value = 48
$ws.Range("H7").Value = 43
$ws.Range("I7").Value = 47

# Update Total row (row 8) values for years 2015-2022 (columns B-I).
$ws.Range("B8").Value = 112
$ws.Range("C8").Value = 225
$ws.Range("D8").Value = 337
$ws.Range("E8").Value = 319
$ws.Range("F8").Value = 218
$ws.Range("G8").Value = 406
$ws.Range("H8").Value = 674
$ws.Range("I8").Value = 710
